$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category") before the existing "date" column,
# shifting date / legislator_name / legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
